# Natmi following Dr Hou advice
# Expand the Cxcl13-Cxcr3 ligand-receptor pair sheet from a single
# FAPs -> FAPs row into the full cross product of sending/target
# clusters (FAPs, M2, sCs) recomputed with the new expression stats.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Cxcl13"
$ws.Range("C2").Value = "Cxcr3"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 4.046465333333333
$ws.Range("H2").Value = 12.139396
$ws.Range("I2").Value = 0.8279246837497715
$ws.Range("J2").Value = 0.8279246837497715
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.05194133333333333
$ws.Range("N2").Value = 0.155824
$ws.Range("O2").Value = 0.03788844568234288
$ws.Range("P2").Value = 0.03788844568234288
$ws.Range("Q2").Value = 0.2101788047004444
$ws.Range("R2").Value = 1.891609242304
$ws.Range("S2").Value = 0.03136877940932412
$ws.Range("T2").Value = 0.03136877940932412
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Cxcl13"
$ws.Range("C3").Value = "Cxcr3"
$ws.Range("D3").Value = "M2"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 4.046465333333333
$ws.Range("H3").Value = 12.139396
$ws.Range("I3").Value = 0.8279246837497715
$ws.Range("J3").Value = 0.8279246837497715
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1.318960333333333
$ws.Range("N3").Value = 3.956881
$ws.Range("O3").Value = 0.962111554317657
$ws.Range("P3").Value = 0.9621115543176572
$ws.Range("Q3").Value = 5.33712726487511
$ws.Range("R3").Value = 48.034145383876
$ws.Range("S3").Value = 0.7965559043404473
$ws.Range("T3").Value = 0.7965559043404473
$ws.Range("A4").Value = "M2"
$ws.Range("B4").Value = "Cxcl13"
$ws.Range("C4").Value = "Cxcr3"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.2754046666666667
$ws.Range("H4").Value = 0.826214
$ws.Range("I4").Value = 0.05634901148785605
$ws.Range("J4").Value = 0.05634901148785604
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.05194133333333333
$ws.Range("N4").Value = 0.155824
$ws.Range("O4").Value = 0.03788844568234288
$ws.Range("P4").Value = 0.03788844568234288
$ws.Range("Q4").Value = 0.01430488559288889
$ws.Range("R4").Value = 0.128743970336
$ws.Range("S4").Value = 0.002134976461011349
$ws.Range("T4").Value = 0.002134976461011349
$ws.Range("A5").Value = "M2"
$ws.Range("B5").Value = "Cxcl13"
$ws.Range("C5").Value = "Cxcr3"
$ws.Range("D5").Value = "M2"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.2754046666666667
$ws.Range("H5").Value = 0.826214
$ws.Range("I5").Value = 0.05634901148785605
$ws.Range("J5").Value = 0.05634901148785604
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.318960333333333
$ws.Range("N5").Value = 3.956881
$ws.Range("O5").Value = 0.962111554317657
$ws.Range("P5").Value = 0.9621115543176572
$ws.Range("Q5").Value = 0.3632478309482222
$ws.Range("R5").Value = 3.269230478534
$ws.Range("S5").Value = 0.0542140350268447
$ws.Range("T5").Value = 0.0542140350268447
$ws.Range("A6").Value = "sCs"
$ws.Range("B6").Value = "Cxcl13"
$ws.Range("C6").Value = "Cxcr3"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.5656100000000001
$ws.Range("H6").Value = 1.69683
$ws.Range("I6").Value = 0.1157263047623724
$ws.Range("J6").Value = 0.1157263047623724
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.05194133333333333
$ws.Range("N6").Value = 0.155824
$ws.Range("O6").Value = 0.03788844568234288
$ws.Range("P6").Value = 0.03788844568234288
$ws.Range("Q6").Value = 0.02937853754666667
$ws.Range("R6").Value = 0.26440683792
$ws.Range("S6").Value = 0.004384689812007407
$ws.Range("T6").Value = 0.004384689812007406
$ws.Range("A7").Value = "sCs"
$ws.Range("B7").Value = "Cxcl13"
$ws.Range("C7").Value = "Cxcr3"
$ws.Range("D7").Value = "M2"
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.5656100000000001
$ws.Range("H7").Value = 1.69683
$ws.Range("I7").Value = 0.1157263047623724
$ws.Range("J7").Value = 0.1157263047623724
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1.318960333333333
$ws.Range("N7").Value = 3.956881
$ws.Range("O7").Value = 0.962111554317657
$ws.Range("P7").Value = 0.9621115543176572
$ws.Range("Q7").Value = 0.7460171541366667
$ws.Range("R7").Value = 6.714154387230001
$ws.Range("S7").Value = 0.111341614950365
$ws.Range("T7").Value = 0.111341614950365